$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 10.18201399862603
$ws.Cells.Item(2, 3).Value = 4.3408564290895
$ws.Cells.Item(2, 4).Value = 8.631474025476173
$ws.Cells.Item(2, 5).Value = 13.63263849530772
$ws.Cells.Item(2, 6).Value = 36.21285343707817
$ws.Cells.Item(2, 8).Value = 7.344005520526261
$ws.Cells.Item(2, 10).Value = 10.32451255953212
$ws.Cells.Item(2, 11).Value = 9.564593265398411
$ws.Cells.Item(2, 13).Value = 15.0264000258118
$ws.Cells.Item(2, 15).Value = 27.63659263531565

$ws.Cells.Item(3, 2).Value = 9.943804821937418
$ws.Cells.Item(3, 3).Value = 4.194357348631359
$ws.Cells.Item(3, 4).Value = 8.59990821162177
$ws.Cells.Item(3, 5).Value = 13.62914680221794
$ws.Cells.Item(3, 6).Value = 36.28459783135221
$ws.Cells.Item(3, 8).Value = 7.344005520526261
$ws.Cells.Item(3, 10).Value = 10.34531993874208
$ws.Cells.Item(3, 11).Value = 9.407726021104821
$ws.Cells.Item(3, 13).Value = 14.96942325103575
$ws.Cells.Item(3, 15).Value = 27.71900279472625

$ws.Cells.Item(4, 2).Value = 9.79648155703935
$ws.Cells.Item(4, 3).Value = 4.101062199758081
$ws.Cells.Item(4, 4).Value = 8.581810633105521
$ws.Cells.Item(4, 5).Value = 13.62927148237973
$ws.Cells.Item(4, 6).Value = 36.33633611243476
$ws.Cells.Item(4, 8).Value = 7.344005520526261
$ws.Cells.Item(4, 10).Value = 10.35926522512722
$ws.Cells.Item(4, 11).Value = 9.311584336855459
$ws.Cells.Item(4, 13).Value = 14.93664562874706
$ws.Cells.Item(4, 15).Value = 27.77486253062047

$ws.Cells.Item(5, 2).Value = 9.736269120299125
$ws.Cells.Item(5, 3).Value = 4.062240146958568
$ws.Cells.Item(5, 4).Value = 8.574763920545768
$ws.Cells.Item(5, 5).Value = 13.62989428276139
$ws.Cells.Item(5, 6).Value = 36.35934970794551
$ws.Cells.Item(5, 8).Value = 7.344005520526261
$ws.Cells.Item(5, 10).Value = 10.36524240525651
$ws.Cells.Item(5, 11).Value = 9.272498714816161
$ws.Cells.Item(5, 13).Value = 14.92385294399508
$ws.Cells.Item(5, 15).Value = 27.79894645412401

$ws.Cells.Item(6, 2).Value = 9.72626289611234
$ws.Cells.Item(6, 3).Value = 4.055746414388293
$ws.Cells.Item(6, 4).Value = 8.573613790983476
$ws.Cells.Item(6, 5).Value = 13.63003228520347
$ws.Cells.Item(6, 6).Value = 36.36328755531137
$ws.Cells.Item(6, 8).Value = 7.344005520526261
$ws.Cells.Item(6, 10).Value = 10.3662526986388
$ws.Cells.Item(6, 11).Value = 9.266015621596258
$ws.Cells.Item(6, 13).Value = 14.92176309133654
$ws.Cells.Item(6, 15).Value = 27.80302527584571

$ws.Cells.Item(7, 2).Value = 9.795670106359005
$ws.Cells.Item(7, 3).Value = 4.10054183309964
$ws.Cells.Item(7, 4).Value = 8.581714263008626
$ws.Cells.Item(7, 5).Value = 13.6292775638601
$ws.Cells.Item(7, 6).Value = 36.33663867249508
$ws.Cells.Item(7, 8).Value = 7.344005520526261
$ws.Cells.Item(7, 10).Value = 10.35934464322507
$ws.Cells.Item(7, 11).Value = 9.311056773591835
$ws.Cells.Item(7, 13).Value = 14.93647080431492
$ws.Cells.Item(7, 15).Value = 27.77518199029946

$ws.Cells.Item(8, 2).Value = 10.10015277305878
$ws.Cells.Item(8, 3).Value = 4.291058831483331
$ws.Cells.Item(8, 4).Value = 8.620326945907603
$ws.Cells.Item(8, 5).Value = 13.63096470493361
$ws.Cells.Item(8, 6).Value = 36.2359937293227
$ws.Cells.Item(8, 8).Value = 7.344005520526261
$ws.Cells.Item(8, 10).Value = 10.33144439581676
$ws.Cells.Item(8, 11).Value = 9.510498130960791
$ws.Cells.Item(8, 13).Value = 15.00630247138898
$ws.Cells.Item(8, 15).Value = 27.66391489844217

$ws.Cells.Item(9, 2).Value = 10.68505201329494
$ws.Cells.Item(9, 3).Value = 4.636682066171032
$ws.Cells.Item(9, 4).Value = 8.705985580189003
$ws.Cells.Item(9, 5).Value = 13.65219697096641
$ws.Cells.Item(9, 6).Value = 36.09974565224474
$ws.Cells.Item(9, 8).Value = 7.344005520526261
$ws.Cells.Item(9, 10).Value = 10.28599846796049
$ws.Cells.Item(9, 11).Value = 9.900883870613864
$ws.Cells.Item(9, 13).Value = 15.16029985778434
$ws.Cells.Item(9, 15).Value = 27.48753658605562

$ws.Cells.Item(10, 2).Value = 11.10260869652089
$ws.Cells.Item(10, 3).Value = 4.871874571669962
$ws.Cells.Item(10, 4).Value = 8.774626225823971
$ws.Cells.Item(10, 5).Value = 13.67860724859336
$ws.Cells.Item(10, 6).Value = 36.03704875132544
$ws.Cells.Item(10, 8).Value = 7.344005520526261
$ws.Cells.Item(10, 10).Value = 10.2582408680765
$ws.Cells.Item(10, 11).Value = 10.18450724213172
$ws.Cells.Item(10, 13).Value = 15.28320591447971
$ws.Cells.Item(10, 15).Value = 27.38355434984815

$ws.Cells.Item(11, 2).Value = 11.28901741444841
$ws.Cells.Item(11, 3).Value = 4.974512900817038
$ws.Cells.Item(11, 4).Value = 8.807010973570224
$ws.Cells.Item(11, 5).Value = 13.69293948522035
$ws.Cells.Item(11, 6).Value = 36.01667000668041
$ws.Cells.Item(11, 8).Value = 7.344005520526261
$ws.Cells.Item(11, 10).Value = 10.2468322349153
$ws.Cells.Item(11, 11).Value = 10.31228527292691
$ws.Cells.Item(11, 13).Value = 15.34109378105425
$ws.Cells.Item(11, 15).Value = 27.34183146054829

$ws.Cells.Item(12, 2).Value = 11.35902794234378
$ws.Cells.Item(12, 3).Value = 5.012733402440536
$ws.Cells.Item(12, 4).Value = 8.819433755927758
$ws.Cells.Item(12, 5).Value = 13.69869705613013
$ws.Cells.Item(12, 6).Value = 36.01012481943762
$ws.Cells.Item(12, 8).Value = 7.344005520526261
$ws.Cells.Item(12, 10).Value = 10.24268697846641
$ws.Cells.Item(12, 11).Value = 10.36044858737735
$ws.Cells.Item(12, 13).Value = 15.36328573881586
$ws.Cells.Item(12, 15).Value = 27.32683601417276

$ws.Cells.Item(13, 2).Value = 11.34397669593837
$ws.Cells.Item(13, 3).Value = 5.004530982570095
$ws.Cells.Item(13, 4).Value = 8.816751319081956
$ws.Cells.Item(13, 5).Value = 13.69744242336727
$ws.Cells.Item(13, 6).Value = 36.01148231738255
$ws.Cells.Item(13, 8).Value = 7.344005520526261
$ws.Cells.Item(13, 10).Value = 10.24357195716825
$ws.Cells.Item(13, 11).Value = 10.35008638274634
$ws.Cells.Item(13, 13).Value = 15.35849446245187
$ws.Cells.Item(13, 15).Value = 27.33002976113052

$ws.Cells.Item(14, 2).Value = 11.29478919350357
$ws.Cells.Item(14, 3).Value = 4.97767038800917
$ws.Cells.Item(14, 4).Value = 8.808029852268898
$ws.Cells.Item(14, 5).Value = 13.69340656359731
$ws.Cells.Item(14, 6).Value = 36.01610804213607
$ws.Cells.Item(14, 8).Value = 7.344005520526261
$ws.Cells.Item(14, 10).Value = 10.24648769755085
$ws.Cells.Item(14, 11).Value = 10.31625242446056
$ws.Cells.Item(14, 13).Value = 15.34291417123236
$ws.Cells.Item(14, 15).Value = 27.34058165060171

$ws.Cells.Item(15, 2).Value = 11.26458313325079
$ws.Cells.Item(15, 3).Value = 4.96113275108816
$ws.Cells.Item(15, 4).Value = 8.802708234843015
$ws.Cells.Item(15, 5).Value = 13.69097739738849
$ws.Cells.Item(15, 6).Value = 36.01909405363033
$ws.Cells.Item(15, 8).Value = 7.344005520526261
$ws.Cells.Item(15, 10).Value = 10.24829644778464
$ws.Cells.Item(15, 11).Value = 10.29549775832889
$ws.Cells.Item(15, 13).Value = 15.33340569552178
$ws.Cells.Item(15, 15).Value = 27.34714975825466

$ws.Cells.Item(16, 2).Value = 11.09034916586432
$ws.Cells.Item(16, 3).Value = 4.865077520724693
$ws.Cells.Item(16, 4).Value = 8.772532517367873
$ws.Cells.Item(16, 5).Value = 13.67771699220091
$ws.Cells.Item(16, 6).Value = 36.03854448831281
$ws.Cells.Item(16, 8).Value = 7.344005520526261
$ws.Cells.Item(16, 10).Value = 10.25901094473577
$ws.Cells.Item(16, 11).Value = 10.17612776209553
$ws.Cells.Item(16, 13).Value = 15.27946147051788
$ws.Cells.Item(16, 15).Value = 27.38639346482502

$ws.Cells.Item(17, 2).Value = 10.98250458335516
$ws.Cells.Item(17, 3).Value = 4.805020591414053
$ws.Cells.Item(17, 4).Value = 8.754312552865448
$ws.Cells.Item(17, 5).Value = 13.6701737752732
$ws.Cells.Item(17, 6).Value = 36.05256299944553
$ws.Cells.Item(17, 8).Value = 7.344005520526261
$ws.Cells.Item(17, 10).Value = 10.26589582702955
$ws.Cells.Item(17, 11).Value = 10.10254723044542
$ws.Cells.Item(17, 13).Value = 15.24686556467198
$ws.Cells.Item(17, 15).Value = 27.41189850271033

$ws.Cells.Item(18, 2).Value = 10.92014585289786
$ws.Cells.Item(18, 3).Value = 4.770068861941536
$ws.Cells.Item(18, 4).Value = 8.743942532575407
$ws.Cells.Item(18, 5).Value = 13.66605350477998
$ws.Cells.Item(18, 6).Value = 36.06139245471559
$ws.Cells.Item(18, 8).Value = 7.344005520526261
$ws.Cells.Item(18, 10).Value = 10.26997053070475
$ws.Cells.Item(18, 11).Value = 10.06011114236024
$ws.Cells.Item(18, 13).Value = 15.22830419283092
$ws.Cells.Item(18, 15).Value = 27.42709335984725

$ws.Cells.Item(19, 2).Value = 10.8989778871836
$ws.Cells.Item(19, 3).Value = 4.758165306472612
$ws.Cells.Item(19, 4).Value = 8.740450470875114
$ws.Cells.Item(19, 5).Value = 13.6646960496291
$ws.Cells.Item(19, 6).Value = 36.06451354542014
$ws.Cells.Item(19, 8).Value = 7.344005520526261
$ws.Cells.Item(19, 10).Value = 10.27136986338278
$ws.Cells.Item(19, 11).Value = 10.04572473652999
$ws.Cells.Item(19, 13).Value = 15.22205213408512
$ws.Cells.Item(19, 15).Value = 27.43232819476689

$ws.Cells.Item(20, 2).Value = 10.99401943199578
$ws.Cells.Item(20, 3).Value = 4.811456197522748
$ws.Cells.Item(20, 4).Value = 8.756240806610434
$ws.Cells.Item(20, 5).Value = 13.67095418262678
$ws.Cells.Item(20, 6).Value = 36.05099137918276
$ws.Cells.Item(20, 8).Value = 7.344005520526261
$ws.Cells.Item(20, 10).Value = 10.26515105027477
$ws.Cells.Item(20, 11).Value = 10.11039218480734
$ws.Cells.Item(20, 13).Value = 15.25031620257186
$ws.Cells.Item(20, 15).Value = 27.40912909518138

$ws.Cells.Item(21, 2).Value = 11.30925297986335
$ws.Cells.Item(21, 3).Value = 4.985577690056044
$ws.Cells.Item(21, 4).Value = 8.810587292969268
$ws.Cells.Item(21, 5).Value = 13.69458305596584
$ws.Cells.Item(21, 6).Value = 36.01471754756026
$ws.Cells.Item(21, 8).Value = 7.344005520526261
$ws.Cells.Item(21, 10).Value = 10.24562652810386
$ws.Cells.Item(21, 11).Value = 10.32619668810456
$ws.Cells.Item(21, 13).Value = 15.34748323096334
$ws.Cells.Item(21, 15).Value = 27.33746046767083

$ws.Cells.Item(22, 2).Value = 11.51187645310356
$ws.Cells.Item(22, 3).Value = 5.095600668218506
$ws.Cells.Item(22, 4).Value = 8.847031274620591
$ws.Cells.Item(22, 5).Value = 13.71194931071961
$ws.Cells.Item(22, 6).Value = 35.99784070020738
$ws.Cells.Item(22, 8).Value = 7.344005520526261
$ws.Cells.Item(22, 10).Value = 10.23388571106598
$ws.Cells.Item(22, 11).Value = 10.46591571661016
$ws.Cells.Item(22, 13).Value = 15.41256114609169
$ws.Cells.Item(22, 15).Value = 27.29530858694944

$ws.Cells.Item(23, 2).Value = 11.40406513968334
$ws.Cells.Item(23, 3).Value = 5.037230922739875
$ws.Cells.Item(23, 4).Value = 8.827498260510385
$ws.Cells.Item(23, 5).Value = 13.7025057036274
$ws.Cells.Item(23, 6).Value = 36.0062230606206
$ws.Cells.Item(23, 8).Value = 7.344005520526261
$ws.Cells.Item(23, 10).Value = 10.24005880328268
$ws.Cells.Item(23, 11).Value = 10.39148018565432
$ws.Cells.Item(23, 13).Value = 15.37768832538335
$ws.Cells.Item(23, 15).Value = 27.31737633956308

$ws.Cells.Item(24, 2).Value = 10.98881467695243
$ws.Cells.Item(24, 3).Value = 4.808547979723423
$ws.Cells.Item(24, 4).Value = 8.755368715856321
$ws.Cells.Item(24, 5).Value = 13.67060068609868
$ws.Cells.Item(24, 6).Value = 36.05169951017774
$ws.Cells.Item(24, 8).Value = 7.344005520526261
$ws.Cells.Item(24, 10).Value = 10.26548740095924
$ws.Cells.Item(24, 11).Value = 10.10684589384438
$ws.Cells.Item(24, 13).Value = 15.24875561232817
$ws.Cells.Item(24, 15).Value = 27.41037948831603

$ws.Cells.Item(25, 2).Value = 10.52863372414443
$ws.Cells.Item(25, 3).Value = 4.546364004596921
$ws.Cells.Item(25, 4).Value = 8.681783727214937
$ws.Cells.Item(25, 5).Value = 13.64454513077241
$ws.Cells.Item(25, 6).Value = 36.13004440825696
$ws.Cells.Item(25, 8).Value = 7.344005520526261
$ws.Cells.Item(25, 10).Value = 10.29730255044693
$ws.Cells.Item(25, 11).Value = 9.795640552352998
$ws.Cells.Item(25, 13).Value = 15.11687678985123
$ws.Cells.Item(25, 15).Value = 27.53076307467229

$wb.Save()